$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin data (prices and volume percentages) scraped on Fri Apr  5 11:30:57 UTC 2024
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.782.72'
$ws.Range('E2').Value = '  +0.98%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.265.71'
$ws.Range('E3').Value = '  -2.19%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.38%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.48'
$ws.Range('E5').Value = '  -1.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.57'
$ws.Range('E6').Value = '  -7.31%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.576'
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.264.81'
$ws.Range('E9').Value = '  -2.10%  '
$ws.Range('E10').Value = '  -5.17%  '
$ws.Range('E11').Value = '  -1.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '44.65'
$ws.Range('E12').Value = '  -4.84%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000267'
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '680.39'
$ws.Range('E14').Value = '  +3.85%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.797.43'
$ws.Range('E15').Value = '  -1.62%  '
$ws.Range('E16').Value = '  -3.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.991.30'
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range('E18').Value = '  +0.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.265.73'
$ws.Range('E19').Value = '  -2.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.12'
$ws.Range('E20').Value = '  -4.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.58'
$ws.Range('E21').Value = '  -4.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.877'
$ws.Range('E22').Value = '  -2.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '16.80'
$ws.Range('E23').Value = '  -4.70%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.19'
$ws.Range('E24').Value = '  +2.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '97.41'
$ws.Range('E25').Value = '  -2.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.83'
$ws.Range('E26').Value = '  -4.51%  '
$ws.Range('E27').Value = '  -5.65%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.90'
$ws.Range('E28').Value = '  -6.78%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '32.44'
$ws.Range('E29').Value = '  +1.24%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.23'
$ws.Range('E30').Value = '  -3.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.59'
$ws.Range('E31').Value = '  -2.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '574.81'
$ws.Range('E32').Value = '  -4.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '10.79'
$ws.Range('E33').Value = '  -2.48%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.790.75'
$ws.Range('E34').Value = '  -2.37%  '
$ws.Range('E35').Value = '  -3.36%  '
$ws.Range('E36').Value = '  -0.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '55.18'
$ws.Range('E37').Value = '  -2.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.27'
$ws.Range('E38').Value = '  -15.48%  '
$ws.Range('E39').Value = '  -1.24%  '
$ws.Range('E40').Value = '  +0.71%  '
$ws.Range('E41').Value = '  -7.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '31.19'
$ws.Range('E42').Value = '  -5.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0₃0650'
$ws.Range('E43').Value = '  -6.84%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.322'
$ws.Range('E44').Value = '  -5.03%  '
$ws.Range('E45').Value = '  -7.26%  '
$ws.Range('E46').Value = '  -4.07%  '
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('E48').Value = '  -1.34%  '
$ws.Range('E49').Value = '  -1.71%  '
$ws.Range('E50').Value = '  +0.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '127.16'
$ws.Range('E51').Value = '  -1.65%  '
